# Auto-generated Excel COM-interop script
# Refreshes the "cryptos" price table on the active sheet, mirroring the
# GitHub Actions data-pull commit "Updated cryptos list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text (e.g. "25.828.63", "0.00001160").
# Assigning a numeric-looking string through COM .Value would normally let
# Excel reinterpret it as a real number (dropping trailing zeros / using
# scientific notation), so for those cells we temporarily force a Text
# number format, assign the literal string, then restore the original
# cell style so formatting stays exactly as it was.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '25.854.99'
$ws.Range("E2").Value = '  +8.41%  '

Set-TextValue $ws.Range("D3") '1.761.61'
$ws.Range("E3").Value = '  +6.16%  '

$ws.Range("E4").Value = '  +0.89%  '

Set-TextValue $ws.Range("D5") '316.89'
$ws.Range("E5").Value = '  +3.33%  '

Set-TextValue $ws.Range("D6") '0.9962'
$ws.Range("E6").Value = '  +1.20%  '

Set-TextValue $ws.Range("D7") '0.3845'
$ws.Range("E7").Value = '  +3.39%  '

$ws.Range("E8").Value = '  +6.03%  '

Set-TextValue $ws.Range("D9") '51.24'
$ws.Range("E9").Value = '  +7.09%  '

Set-TextValue $ws.Range("D10") '1.239'
$ws.Range("E10").Value = '  +6.19%  '

Set-TextValue $ws.Range("D11") '0.07665'
$ws.Range("E11").Value = '  +6.49%  '

Set-TextValue $ws.Range("D12") '0.9981'
$ws.Range("E12").Value = '  +0.82%  '

Set-TextValue $ws.Range("D13") '21.85'
$ws.Range("E13").Value = '  +6.07%  '

Set-TextValue $ws.Range("D14") '6.494'
$ws.Range("E14").Value = '  +8.42%  '

Set-TextValue $ws.Range("D15") '7.147'
$ws.Range("E15").Value = '  +6.06%  '

Set-TextValue $ws.Range("D16") '1.765.62'
$ws.Range("E16").Value = '  +6.11%  '

Set-TextValue $ws.Range("D17") '0.00001163'
$ws.Range("E17").Value = '  +6.34%  '

Set-TextValue $ws.Range("D18") '0.9966'
$ws.Range("E18").Value = '  +1.32%  '

Set-TextValue $ws.Range("D19") '0.06859'
$ws.Range("E19").Value = '  +1.94%  '

Set-TextValue $ws.Range("D20") '87.44'
$ws.Range("E20").Value = '  +7.73%  '

Set-TextValue $ws.Range("D21") '17.78'
$ws.Range("E21").Value = '  +8.48%  '

Set-TextValue $ws.Range("D22") '6.555'
$ws.Range("E22").Value = '  +7.97%  '

Set-TextValue $ws.Range("D23") '12.82'
$ws.Range("E23").Value = '  +7.44%  '

Set-TextValue $ws.Range("D24") '25.837.25'
$ws.Range("E24").Value = '  +8.39%  '

Set-TextValue $ws.Range("D25") '2.434'
$ws.Range("E25").Value = '  +3.72%  '

Set-TextValue $ws.Range("D26") '2.984'
$ws.Range("E26").Value = '  +11.34%  '

Set-TextValue $ws.Range("D27") '20.80'
$ws.Range("E27").Value = '  +6.66%  '

Set-TextValue $ws.Range("D28") '155.11'
$ws.Range("E28").Value = '  +2.22%  '

$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D29") '135.49'
$ws.Range("E29").Value = '  +7.16%  '

$ws.Range("B30").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C30").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D30") '1.963.54'
$ws.Range("E30").Value = '  +6.24%  '

Set-TextValue $ws.Range("D31") '1.202'
$ws.Range("E31").Value = '  +21.88%  '

Set-TextValue $ws.Range("D32") '7.233'
$ws.Range("E32").Value = '  +18.39%  '

Set-TextValue $ws.Range("D33") '4.320'
$ws.Range("E33").Value = '  +7.26%  '

Set-TextValue $ws.Range("D34") '14.05'
$ws.Range("E34").Value = '  +15.20%  '

Set-TextValue $ws.Range("D35") '1.809'
$ws.Range("E35").Value = '  +6.22%  '

Set-TextValue $ws.Range("D36") '0.08729'
$ws.Range("E36").Value = '  +4.35%  '

Set-TextValue $ws.Range("D37") '5.693'
$ws.Range("E37").Value = '  +8.13%  '

Set-TextValue $ws.Range("D38") '0.06784'
$ws.Range("E38").Value = '  +7.42%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D39") '0.02492'
$ws.Range("E39").Value = '  +9.47%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D40") '9.333'
$ws.Range("E40").Value = '  +5.52%  '

Set-TextValue $ws.Range("D41") '0.2241'
$ws.Range("E41").Value = '  +8.90%  '

Set-TextValue $ws.Range("D42") '1.308'
$ws.Range("E42").Value = '  +3.46%  '

Set-TextValue $ws.Range("D43") '0.6600'
$ws.Range("E43").Value = '  +9.45%  '

Set-TextValue $ws.Range("D44") '14.27'
$ws.Range("E44").Value = '  +7.72%  '

Set-TextValue $ws.Range("D45") '0.9964'
$ws.Range("E45").Value = '  +1.66%  '

Set-TextValue $ws.Range("D46") '0.6380'
$ws.Range("E46").Value = '  +8.32%  '

Set-TextValue $ws.Range("D47") '3.921'
$ws.Range("E47").Value = '  +2.41%  '

Set-TextValue $ws.Range("D48") '2.179'
$ws.Range("E48").Value = '  +9.34%  '

Set-TextValue $ws.Range("D49") '133.65'
$ws.Range("E49").Value = '  +5.48%  '

Set-TextValue $ws.Range("D50") '0.07505'
$ws.Range("E50").Value = '  +6.19%  '

Set-TextValue $ws.Range("D51") '81.03'
$ws.Range("E51").Value = '  +7.43%  '

